$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so COM does not
# reinterpret numeric-looking strings (e.g. "232.41") as numbers,
# then restore the original (default) style once values are set.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.330.38"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.353.04"
$ws.Range("E3").Value = "  +5.50%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "232.41"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "0.648"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").Value = "67.36"
$ws.Range("E7").Value = "  +6.91%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").Value = "0.0960"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "56.81"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "26.30"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "2.701.90"
$ws.Range("E13").Value = "  +5.40%  "
$ws.Range("D15").Value = "15.69"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "2.347.29"
$ws.Range("E18").Value = "  +5.01%  "
$ws.Range("D19").Value = "43.251.32"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "73.87"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("D23").Value = "248.78"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "3.97"
$ws.Range("E24").Value = "  +17.24%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "22.36"
$ws.Range("E29").Value = "  +7.63%  "
$ws.Range("D30").Value = "171.95"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "1.53"
$ws.Range("E31").Value = "  +11.55%  "
$ws.Range("E32").Value = "  -7.48%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "5.03"
$ws.Range("E34").Value = "  +6.46%  "
$ws.Range("D35").Value = "0.0691"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("E37").Value = "  +10.23%  "
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").Value = "3.61"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("E41").Value = "  +8.86%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "18.23"
$ws.Range("E43").Value = "  +6.90%  "
$ws.Range("E44").Value = "  +9.02%  "
$ws.Range("D45").Value = "1.21"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "4.47"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").Value = "98.30"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "0.0951"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "1.443.98"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").Value = "2.574.68"
$ws.Range("E50").Value = "  +5.62%  "
$ws.Range("E51").Value = "  -1.73%  "

$dataRange.Style = "Normal"
